$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.0743280361453236
$ws.Range("H2").Value = 14.78948713255983
$ws.Range("I2").Value = 475.1119739588566
$ws.Range("G3").Value = 0.06825074401602975
$ws.Range("H3").Value = -0.2256082236109598
$ws.Range("G4").Value = -0.03915086399347909
$ws.Range("H4").Value = 12.79078838075259
$ws.Range("G5").Value = -0.01866125984821869
$ws.Range("H5").Value = 31.0676403719955
$ws.Range("G6").Value = -0.08271557824885563
$ws.Range("H6").Value = 21.99675262930482
$ws.Range("G7").Value = -0.09015388193962762
$ws.Range("H7").Value = 1.332242977728449
$ws.Range("G8").Value = -0.3599841737535226
$ws.Range("H8").Value = 1.886146153813329
$ws.Range("G9").Value = -0.3648216446658841
$ws.Range("H9").Value = 6.477744396467491
$ws.Range("G10").Value = 0.03807514592041499
$ws.Range("H10").Value = 88.57485923981226
$ws.Range("G11").Value = 0.02842654780157124
$ws.Range("H11").Value = 25.25756375653284
$ws.Range("G12").Value = 0.2142745085071594
$ws.Range("H12").Value = -3.370067637944354
$ws.Range("G13").Value = 0.2241120825662148
$ws.Range("H13").Value = -0.4860949158006961
$ws.Range("G14").Value = -0.04543580366458032
$ws.Range("H14").Value = -7.910444908633031
$ws.Range("G15").Value = -0.05472133593534471
$ws.Range("H15").Value = -14.72868958465596
$ws.Range("G16").Value = 0.216025536793977
$ws.Range("H16").Value = 1.62397460901277
$ws.Range("G17").Value = 0.2206189229254771
$ws.Range("H17").Value = 0.03217373801147629
$ws.Range("G18").Value = 0.07658885007158407
$ws.Range("H18").Value = 4.887151084722289
$ws.Range("G19").Value = 0.08043953820566083
$ws.Range("H19").Value = 6.77151577036186
$ws.Range("G20").Value = -0.08044438584332532
$ws.Range("H20").Value = -7.284963689751561
$ws.Range("G21").Value = -0.08234281939936004
$ws.Range("H21").Value = 4.881678450407706
$ws.Range("G22").Value = 0.07323269483254494
$ws.Range("H22").Value = -0.3662517747011098
$ws.Range("G23").Value = 0.07384821209514221
$ws.Range("H23").Value = 8.073744274330243
$ws.Range("G24").Value = 0.06746233037103917
$ws.Range("H24").Value = 1.276342860401429
$ws.Range("G25").Value = 0.06412377492086697
$ws.Range("H25").Value = 17.05679403628633
$ws.Range("G26").Value = 0.1182826026302501
$ws.Range("H26").Value = -0.8935401931388353
$ws.Range("G27").Value = 0.1225704598949948
$ws.Range("H27").Value = 7.649975916714906
$ws.Range("G28").Value = 0.1341822115967917
$ws.Range("H28").Value = 3.813310518706193
$ws.Range("G29").Value = 0.1375336028696899
$ws.Range("H29").Value = -8.821780564892286
$ws.Range("G30").Value = 0.08898821025144922
$ws.Range("H30").Value = 5.552191736007843
$ws.Range("G31").Value = 0.09160242339596379
$ws.Range("H31").Value = 12.13687841450907
$ws.Range("G32").Value = 0.05644785683990616
$ws.Range("H32").Value = 5.787042167804707
$ws.Range("G33").Value = 0.05639984961599955
$ws.Range("H33").Value = 2.094275561666261
$ws.Range("G34").Value = 0.01013449298524446
$ws.Range("H34").Value = -41.61231518947289
$ws.Range("G35").Value = 0.02139535201274146
$ws.Range("H35").Value = 26.59456135812964
$ws.Range("G36").Value = -0.0204480150744816
$ws.Range("H36").Value = 29.60054423824313
$ws.Range("G37").Value = -0.01589760301219547
$ws.Range("H37").Value = 52.20760883615192
$ws.Range("G38").Value = 0.07632324828959812
$ws.Range("H38").Value = -2.503074278795301
$ws.Range("G39").Value = 0.07945846834961998
$ws.Range("H39").Value = 2.203787428844034
$ws.Range("G40").Value = 0.06392746202049054
$ws.Range("H40").Value = -3.446009631985637
$ws.Range("G41").Value = 0.06509649089041833
$ws.Range("H41").Value = 0.1132700999495589
$ws.Range("G42").Value = 0.07744504374521131
$ws.Range("H42").Value = -0.4451004329846407
$ws.Range("G43").Value = 0.08768242427589364
$ws.Range("H43").Value = 9.378132434410915
$ws.Range("G44").Value = 0.07851875215850995
$ws.Range("H44").Value = -11.02418963948182
$ws.Range("G45").Value = 0.09625700292452216
$ws.Range("H45").Value = 6.498417237172768
$ws.Range("G46").Value = 0.001106459218782339
$ws.Range("H46").Value = 140.4374882123752
$ws.Range("G47").Value = 0.01244202760357919
$ws.Range("H47").Value = 26143.84149518692
$ws.Range("G48").Value = -0.08807072099565647
$ws.Range("H48").Value = 8.360034550299268
$ws.Range("G49").Value = -0.102486770484272
$ws.Range("H49").Value = 6.469865186568917
$ws.Range("G50").Value = 0.1778830459905756
$ws.Range("H50").Value = 4.330666130317605
$ws.Range("G51").Value = 0.1726644126136901
$ws.Range("H51").Value = 1.666055983233625
$ws.Range("G52").Value = 0.06226014594962717
$ws.Range("H52").Value = -12.26255290692823
$ws.Range("G53").Value = 0.07164401204000689
$ws.Range("H53").Value = 11.39357489440809
$ws.Range("G54").Value = -0.1336095592000882
$ws.Range("H54").Value = -4.543484979302187
$ws.Range("G55").Value = -0.1224624651157438
$ws.Range("H55").Value = -5.144017743809016
$ws.Range("G56").Value = 0.1901323320766818
$ws.Range("H56").Value = 0.05867617178333032
$ws.Range("G57").Value = 0.2036309055125753
$ws.Range("H57").Value = 2.374275865142782
